$wb = $excel.ActiveWorkbook

# Sheet ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 90.818184
$ws.Range("I11").Value = 90.818184
$ws.Range("K11").Value = 90.818184
$ws.Range("M11").Value = 49.181816

# Sheet ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1454.1
$ws.Range("I19").Value = 1604.8334
$ws.Range("K19").Value = 1604.8334
$ws.Range("M19").Value = -1429.8334

# Sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 142
$ws.Range("I33").Value = 149
$ws.Range("K33").Value = 149
$ws.Range("M33").Value = 80

# Sheet ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 300
$ws.Range("I53").Value = 100
$ws.Range("K53").Value = 100
$ws.Range("M53").Value = 537

# Sheet ALC row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 772.625
$ws.Range("I55").Value = 777.6
$ws.Range("J55").Value = 764.3333
$ws.Range("K55").Value = 777.6
$ws.Range("L55").Value = 764.3333
$ws.Range("M55").Value = -563.6
$ws.Range("N55").Value = -1192.3333

# Sheet ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1385.4
$ws.Range("I80").Value = 710.4
$ws.Range("J80").Value = 2060.4
$ws.Range("K80").Value = 2131.2
$ws.Range("L80").Value = 6181.200000000001
$ws.Range("M80").Value = -1133.2
$ws.Range("N80").Value = -8177.200000000001

# Sheet ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1385.4
$ws.Range("I83").Value = 710.4
$ws.Range("J83").Value = 2060.4
$ws.Range("K83").Value = 6393.599999999999
$ws.Range("L83").Value = 18543.6
$ws.Range("M83").Value = -1401.599999999999
$ws.Range("N83").Value = -28527.6

# Sheet ARM row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 10226
$ws.Range("I62").Value = 10226
$ws.Range("K62").Value = 10226
$ws.Range("M62").Value = -9602

# Sheet ARM row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 10226
$ws.Range("I65").Value = 10226
$ws.Range("K65").Value = 30678
$ws.Range("M65").Value = -27558

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5263.727
$ws.Range("I132").Value = 560.1429000000001
$ws.Range("J132").Value = 13495
$ws.Range("K132").Value = 1680.4287
$ws.Range("L132").Value = 40485
$ws.Range("M132").Value = 849.5712999999998
$ws.Range("N132").Value = -45545

# Sheet BSM row 11
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 900
$ws.Range("J11").Value = 900
$ws.Range("L11").Value = 900
$ws.Range("N11").Value = -1180

# Sheet CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 444
$ws.Range("I22").Value = 221.8
$ws.Range("K22").Value = 221.8
$ws.Range("M22").Value = 128.2

# Sheet CRP row 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 44995
$ws.Range("J48").Value = 44995
$ws.Range("L48").Value = 44995
$ws.Range("N48").Value = -45947

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Sheet CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 100000
$ws.Range("I18").Value = 100000
$ws.Range("K18").Value = 300000
$ws.Range("M18").Value = -299831

# Sheet CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 119.25
$ws.Range("I23").Value = 152
$ws.Range("K23").Value = 456
$ws.Range("M23").Value = -221

# Sheet CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 594.6667
$ws.Range("I33").Value = 500.75
$ws.Range("K33").Value = 3004.5
$ws.Range("M33").Value = -2721.5

# Sheet CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 979.8570999999999
$ws.Range("J114").Value = 1029
$ws.Range("L114").Value = 3087
$ws.Range("N114").Value = -9595

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1666.6666
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

# Sheet LTW row 10
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5725
$ws.Range("I46").Value = 7900
$ws.Range("K46").Value = 7900
$ws.Range("M46").Value = -7712

# Sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 750.0909
$ws.Range("I55").Value = 583.44446
$ws.Range("K55").Value = 583.44446
$ws.Range("M55").Value = -410.44446

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22500
$ws.Range("J132").Value = 24500
$ws.Range("L132").Value = 73500
$ws.Range("N132").Value = -78560

# Sheet WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4499
$ws.Range("I2").Value = 4499
$ws.Range("K2").Value = 4499
$ws.Range("M2").Value = -4387

# Sheet WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4998.5
$ws.Range("I4").Value = 4998.5
$ws.Range("K4").Value = 4998.5
$ws.Range("M4").Value = -4885.5

# Sheet WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29999.834
$ws.Range("I54").Value = 29999.834
$ws.Range("K54").Value = 29999.834
$ws.Range("M54").Value = -29479.834

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7284.2856
$ws.Range("I81").Value = 390
$ws.Range("J81").Value = 8433.333000000001
$ws.Range("K81").Value = 780
$ws.Range("L81").Value = 16866.666
$ws.Range("M81").Value = 281
$ws.Range("N81").Value = -18988.666

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7284.2856
$ws.Range("I84").Value = 390
$ws.Range("J84").Value = 8433.333000000001
$ws.Range("K84").Value = 3900
$ws.Range("L84").Value = 84333.33
$ws.Range("M84").Value = 1404
$ws.Range("N84").Value = -94941.33

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5501.7896
$ws.Range("I132").Value = 3249.2307
$ws.Range("K132").Value = 9747.6921
$ws.Range("M132").Value = -7217.6921

Write-Output "All edits applied"